$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "NCT(2.118473494661721, 1.281438240921506, -0.2153791949344146, 2.204856187604496)"
$ws.Range("C2").Value = "NIG(0.8757519013663837, 0.5269513092734076, 6.580277629810109, 5.3201469237803085)"
$ws.Range("D2").Value = "JSU(-0.9494824841749, 1.2304715573297296, 0.6480674894796262, 3.0577799086488846)"
$ws.Range("E2").Value = "NCT(2.1350229007052595, 2.13727359457783, -0.011531013113138856, 3.3343369180954863)"
